$d = $word.ActiveDocument

# Helper: replace the ENTIRE contents of a paragraph (by 1-based Paragraphs index)
# with a fresh set of run(s)/content described by $innerXml. Using InsertXML on a
# Range that spans the whole paragraph (including its end-of-paragraph mark)
# replaces the paragraph's content in place while preserving however many
# separate <w:r> elements are present in $innerXml (no run-coalescing), which a
# plain Range.Text assignment would not do.
function Replace-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$full.InsertXML($pkg)
}

# Helper: insert fresh run(s) at a collapsed point inside an existing paragraph
# (used to prepend content before the bookmark in the final paragraph, without
# touching that paragraph's own pPr / bookmark elements).
function Insert-RunsXmlAt($pos, $innerXml) {
    $ins = $d.Range($pos, $pos)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$ins.InsertXML($pkg)
}

$rpr = '<w:rPr><w:lang w:val="en-GB"/></w:rPr>'

# --- Paragraph 6: "All our personas ..." becomes the single-run
#     "After some tests with our client (K. Hindriks) ..." paragraph.
$text1 = "After some tests with our client (K. Hindriks) he said us that he had a hard time during the test to see colours of a button. Therefor we are still considering a good and practical solution to this problem."
$inner1 = '<w:pPr>' + $rpr + '</w:pPr>' +
          '<w:r>' + $rpr + '<w:t>' + $text1 + '</w:t></w:r>'
Replace-ParagraphXml 6 $inner1

# --- Paragraph 7: "After some tests ..." (previously split around "Hindriks"
#     with proofErr spell-check marks) becomes the "All our personas ..."
#     paragraph, re-created with its original six-run split.
$inner2 = '<w:pPr>' + $rpr + '</w:pPr>' +
  '<w:r>' + $rpr + '<w:t>All our personas fall in the same category: they all know how to use goal and study artificial intelligence. If we would change the personas to people who would not do any research</w:t></w:r>' +
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> in</w:t></w:r>' +
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> or study</w:t></w:r>' +
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> artificial intelligence,</w:t></w:r>' +
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r>' + $rpr + '<w:t>we would probably get back to the &#8220;game&#8221; as described for children.</w:t></w:r>'
Replace-ParagraphXml 7 $inner2

# --- Final paragraph (holds the _GoBack bookmark): prepend two new runs with
#     the batchrunner / map editor text, ahead of the existing bookmark marks.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$insPos = $last.Range.Start
$inner3 = '<w:r>' + $rpr + '<w:t xml:space="preserve">If we had not made a batchrunner, a small script that runs the program x times in a row on the background, it would have taken researchers huge amounts of time to run the program 1000 times. Not to mention putting everything together and analyse it afterwards. Without the map editor, users could only select the amount of rows and columns (which would turn out to be the amount of </w:t></w:r>' +
           '<w:r>' + $rpr + '<w:t xml:space="preserve">rooms). The user would not be able to really specify which rooms would be placed on what places and where the drop zone would be. </w:t></w:r>'
Insert-RunsXmlAt $insPos $inner3
